$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF holds the "Date" values that were recorded as text "5-4-2007-08"
# but actually represent 2008-05-04 (NBA stats for the date were shown one
# day off). Fix rows 2 through 31 in column BF (58), keeping the corrected
# value as plain text (not an Excel date serial).
for ($r = 2; $r -le 31; $r++) {
    $cell = $ws.Cells.Item($r, 58)
    if ($cell.Value2 -eq "5-4-2007-08") {
        # Writing the formula for a text literal keeps Excel from
        # re-interpreting "2008-05-04" as a date, then PasteSpecial values
        # collapses it back down to a plain inline string cell (no formula,
        # no style change) - matching a straight text overwrite.
        $cell.Formula = '="2008-05-04"'
        $cell.Copy() | Out-Null
        $cell.PasteSpecial(-4163) | Out-Null
    }
}
$excel.CutCopyMode = $false
